$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so values like
# "57.022.10" or "545.26" are not coerced into numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '57.022.10'
$ws.Range('E2').Value = '  -8.25%  '
$ws.Range('D3').Value = '2.859.16'
$ws.Range('E3').Value = '  -7.54%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '545.26'
$ws.Range('E5').Value = '  -7.69%  '
$ws.Range('D6').Value = '120.42'
$ws.Range('E6').Value = '  -8.54%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '2.852.18'
$ws.Range('E8').Value = '  -7.71%  '
$ws.Range('D9').Value = '0.485'
$ws.Range('E9').Value = '  -3.41%  '
$ws.Range('D10').Value = '0.124'
$ws.Range('E10').Value = '  -11.28%  '
$ws.Range('D11').Value = '4.73'
$ws.Range('E11').Value = '  -10.91%  '
$ws.Range('D12').Value = '0.424'
$ws.Range('E12').Value = '  -4.07%  '
$ws.Range('D13').Value = '0.0000209'
$ws.Range('E13').Value = '  -11.42%  '
$ws.Range('D14').Value = '30.89'
$ws.Range('E14').Value = '  -8.50%  '
$ws.Range('D15').Value = '0.118'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').Value = '3.328.13'
$ws.Range('E16').Value = '  -7.58%  '
$ws.Range('D17').Value = '2.852.86'
$ws.Range('E17').Value = '  -8.06%  '
$ws.Range('D18').Value = '57.064.63'
$ws.Range('E18').Value = '  -8.60%  '
$ws.Range('D19').Value = '6.24'
$ws.Range('E19').Value = '  -2.68%  '
$ws.Range('D20').Value = '407.52'
$ws.Range('E20').Value = '  -9.53%  '
$ws.Range('D21').Value = '12.59'
$ws.Range('E21').Value = '  -7.61%  '
$ws.Range('D22').Value = '0.641'
$ws.Range('E22').Value = '  -5.37%  '
$ws.Range('D23').Value = '6.70'
$ws.Range('E23').Value = '  -9.63%  '
$ws.Range('E24').Value = '  -4.65%  '
$ws.Range('D25').Value = '75.79'
$ws.Range('E25').Value = '  -6.58%  '
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  -7.71%  '
$ws.Range('D29').Value = '1.88'
$ws.Range('E29').Value = '  -7.22%  '
$ws.Range('D30').Value = '7.00'
$ws.Range('E30').Value = '  -6.70%  '
$ws.Range('D31').Value = '24.45'
$ws.Range('E31').Value = '  -6.54%  '
$ws.Range('D32').Value = '5.87'
$ws.Range('E32').Value = '  -9.99%  '
$ws.Range('D33').Value = '0.0918'
$ws.Range('E33').Value = '  -6.53%  '
$ws.Range('D34').Value = '5.29'
$ws.Range('E34').Value = '  -7.39%  '
$ws.Range('D35').Value = '48.16'
$ws.Range('E35').Value = '  -4.80%  '
$ws.Range('D36').Value = '0.879'
$ws.Range('E36').Value = '  -11.37%  '
$ws.Range('D37').Value = '1.96'
$ws.Range('E37').Value = '  -16.38%  '
$ws.Range('D38').Value = '8.15'
$ws.Range('E38').Value = '  +2.52%  '
$ws.Range('D39').Value = '0.0₃0607'
$ws.Range('E39').Value = '  -13.20%  '
$ws.Range('D40').Value = '0.0336'
$ws.Range('E40').Value = '  -11.62%  '
$ws.Range('D41').Value = '0.103'
$ws.Range('E41').Value = '  -6.74%  '
$ws.Range('D42').Value = '2.572.98'
$ws.Range('E42').Value = '  -5.40%  '
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').Value = '347.50'
$ws.Range('E44').Value = '  -8.75%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '2.32'
$ws.Range('E45').Value = '  -8.85%  '
$ws.Range('D46').Value = '116.59'
$ws.Range('E46').Value = '  -7.27%  '
$ws.Range('D47').Value = '0.224'
$ws.Range('E47').Value = '  -7.64%  '
$ws.Range('E48').Value = '  -4.98%  '
$ws.Range('D49').Value = '1.89'
$ws.Range('E49').Value = '  -7.31%  '
$ws.Range('D50').Value = '22.22'
$ws.Range('E50').Value = '  -8.39%  '
$ws.Range('D51').Value = '1.90'
$ws.Range('E51').Value = '  -9.44%  '

# Restore the default (unstyled) cell style on the Price column so the
# text-forcing NumberFormat tweak above does not leave a visible style delta.
$ws.Range('D2:D51').Style = 'Normal'

